# Renamed all output tables .csv's, and completed natl origin table.
# This script adds the new "Natl_Origin" worksheet (with the National Origin
# data table) and registers it in the "TOC" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Update the TOC sheet with a new row describing the new tab ---
$toc = $wb.Worksheets.Item("TOC")
$tocLastRow = $toc.Cells.Item($toc.Rows.Count, 1).End(-4162).Row
$newTocRow = $tocLastRow + 1
$toc.Cells.Item($newTocRow, 1).Value = "Natl_Origin"
$toc.Cells.Item($newTocRow, 2).Value = "National Origin (% Foreign born) by County and SCAG Region"

# --- 2. Add the new "Natl_Origin" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Natl_Origin"

# Header row (bold + centered, matching the other data tables)
$ws.Range("A1").Value = "county"
$ws.Range("B1").Value = "natl_origin_perc"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108

# Data rows
$data = @(
    @("SCAG", 29.32),
    @("Imperial", 29.76),
    @("Los Angeles", 33.51),
    @("Orange", 29.86),
    @("Riverside", 21.55),
    @("San Bernardino", 20.92),
    @("Ventura", 21.2)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
